# Daily attendance processing - 2025-12-31 09:35:11
# Normalizes the "Recorded By" (column G) values so that:
#   - any lowercase "system" token is moved to the end of the list
#   - the canonical "System" token (exact case) is moved to the front
# The rest of the comma-separated tokens keep their relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1
$firstCol = $usedRange.Column
$lastCol = $firstCol + $usedRange.Columns.Count - 1

# Locate the "Recorded By" column from the header row rather than
# hard-coding it, so the script is resilient to column reordering.
$recordedByCol = 7
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $headerText = $ws.Cells.Item($firstRow, $c).Text
    if ($headerText -eq "Recorded By") {
        $recordedByCol = $c
    }
}

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $value = $cell.Text

    if ($value -eq $null) { continue }
    if (-not ($value -is [string])) { continue }
    if ($value.IndexOf(",") -lt 0) { continue }

    $parts = $value.Split(",") | ForEach-Object { $_.Trim() }

    # Case-sensitive partition: exact lowercase "system" tokens move to the end.
    $nonLowerSystem = @($parts | Where-Object { -not $_.Equals("system") })
    $lowerSystem = @($parts | Where-Object { $_.Equals("system") })
    $reordered = @($nonLowerSystem) + @($lowerSystem)

    # Case-sensitive check: exact "System" token moves to the front.
    $hasCanonicalSystem = $false
    foreach ($p in $reordered) {
        if ($p.Equals("System")) { $hasCanonicalSystem = $true }
    }

    if ($hasCanonicalSystem) {
        $rest = @($reordered | Where-Object { -not $_.Equals("System") })
        $final = @("System") + @($rest)
    } else {
        $final = $reordered
    }

    $newValue = [string]::Join(", ", $final)

    if ($newValue -ne $value) {
        $cell.Value = $newValue
    }
}
